# upload TEST CASE UI 2
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UI Test")

# Replace row 5 (previously stray "SES-02" data) with the new "UI-02" form
# validation test case.
$ws.Range("C5").Value = "UI-02"
$ws.Range("E5").Value = "UI-TC-02"
$ws.Range("G5").Value = "UI-TS-02"
$ws.Range("D5").Value = "Verify form validation appears when input is incorrect"
$ws.Range("F5").Value = "UI Module"
$ws.Range("I5").Value = "User is on form page"
$ws.Range("H5").Value = "Verify validation message for`n incorrect input"
$ws.Range("J5").Value = "1. Enter invalid email format `n2. Leave required field empty `n3. Click Submit"
$ws.Range("K5").Value = "Invalid email format, `nempty required field"

$ws.Range("H5").WrapText = $true
$ws.Range("K5").WrapText = $true

$ws.Range("K5").Select()
